$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Stage every changed cell as a formula that evaluates to the literal text
# we need (this avoids Excel auto-converting numeric-looking strings such as
# "1.00" into real numbers when assigned directly via .Value).
$ws.Range("D2").Formula = '="69.859.81"'
$ws.Range("E2").Formula = '="  -0.91%  "'
$ws.Range("D3").Formula = '="3.727.09"'
$ws.Range("E3").Formula = '="  -1.61%  "'
$ws.Range("D4").Formula = '="1.00"'
$ws.Range("E4").Formula = '="  -0.24%  "'
$ws.Range("D5").Formula = '="617.12"'
$ws.Range("E5").Formula = '="  -0.05%  "'
$ws.Range("D6").Formula = '="180.88"'
$ws.Range("E6").Formula = '="  +1.53%  "'
$ws.Range("D7").Formula = '="3.730.07"'
$ws.Range("E7").Formula = '="  -1.14%  "'
$ws.Range("E8").Formula = '="  -0.20%  "'
$ws.Range("E9").Formula = '="  -3.02%  "'
$ws.Range("E10").Formula = '="  -3.35%  "'
$ws.Range("D11").Formula = '="6.33"'
$ws.Range("E11").Formula = '="  -1.22%  "'
$ws.Range("E12").Formula = '="  -4.98%  "'
$ws.Range("D13").Formula = '="39.91"'
$ws.Range("E13").Formula = '="  -1.94%  "'
$ws.Range("E14").Formula = '="  -3.48%  "'
$ws.Range("D15").Formula = '="4.344.03"'
$ws.Range("E15").Formula = '="  -1.96%  "'
$ws.Range("D16").Formula = '="3.726.81"'
$ws.Range("E16").Formula = '="  -2.17%  "'
$ws.Range("D17").Formula = '="69.859.35"'
$ws.Range("E17").Formula = '="  -1.23%  "'
$ws.Range("E18").Formula = '="  -2.11%  "'
$ws.Range("D19").Formula = '="7.55"'
$ws.Range("E19").Formula = '="  -0.84%  "'
$ws.Range("D20").Formula = '="501.16"'
$ws.Range("E20").Formula = '="  -4.53%  "'
$ws.Range("D21").Formula = '="16.30"'
$ws.Range("E21").Formula = '="  -4.05%  "'
$ws.Range("D22").Formula = '="9.29"'
$ws.Range("E22").Formula = '="  -1.36%  "'
$ws.Range("D23").Formula = '="0.720"'
$ws.Range("E23").Formula = '="  -3.65%  "'
$ws.Range("D24").Formula = '="2.53"'
$ws.Range("E24").Formula = '="  +1.55%  "'
$ws.Range("D25").Formula = '="86.42"'
$ws.Range("E25").Formula = '="  -1.89%  "'
$ws.Range("D26").Formula = '="12.94"'
$ws.Range("E26").Formula = '="  -4.68%  "'
$ws.Range("D27").Formula = '="11.19"'
$ws.Range("E27").Formula = '="  +2.59%  "'
$ws.Range("E28").Formula = '="  +5.43%  "'
$ws.Range("D29").Formula = '="1.00"'
$ws.Range("E29").Formula = '="  +0.14%  "'
$ws.Range("E30").Formula = '="  -2.22%  "'
$ws.Range("E31").Formula = '="  -0.26%  "'
$ws.Range("D32").Formula = '="7.94"'
$ws.Range("E32").Formula = '="  -0.25%  "'
$ws.Range("D33").Formula = '="30.31"'
$ws.Range("E33").Formula = '="  -6.27%  "'
$ws.Range("D34").Formula = '="0.114"'
$ws.Range("E34").Formula = '="  -1.23%  "'
$ws.Range("D35").Formula = '="0.999"'
$ws.Range("E35").Formula = '="  -0.34%  "'
$ws.Range("E36").Formula = '="  +0.14%  "'
$ws.Range("E37").Formula = '="  -1.45%  "'
$ws.Range("E38").Formula = '="  +4.25%  "'
$ws.Range("D39").Formula = '="0.344"'
$ws.Range("E39").Formula = '="  +0.01%  "'
$ws.Range("D40").Formula = '="3.08"'
$ws.Range("E40").Formula = '="  +10.43%  "'
$ws.Range("E41").Formula = '="  -6.13%  "'
$ws.Range("D42").Formula = '="50.02"'
$ws.Range("E42").Formula = '="  -3.20%  "'
$ws.Range("D43").Formula = '="425.75"'
$ws.Range("E43").Formula = '="  -1.09%  "'
$ws.Range("D44").Formula = '="44.10"'
$ws.Range("E44").Formula = '="  -0.71%  "'
$ws.Range("E45").Formula = '="  -3.55%  "'
$ws.Range("D46").Formula = '="2.948.63"'
$ws.Range("E46").Formula = '="  -6.57%  "'
$ws.Range("E47").Formula = '="  -2.49%  "'
$ws.Range("B48").Formula = '="USDe"'
$ws.Range("C48").Formula = '="https://coinranking.com/coin/exbfr2U-0+usde-usde"'
$ws.Range("D48").Formula = '="1.00"'
$ws.Range("E48").Formula = '="  -0.08%  "'
$ws.Range("B49").Formula = '="InjectiveProtocol"'
$ws.Range("C49").Formula = '="https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"'
$ws.Range("D49").Formula = '="27.10"'
$ws.Range("E49").Formula = '="  -2.34%  "'
$ws.Range("D50").Formula = '="136.68"'
$ws.Range("E50").Formula = '="  -2.80%  "'
$ws.Range("E51").Formula = '="  -2.93%  "'

# Convert every staged formula back into a plain (static) text value via
# copy / paste-special-values, exactly as "Paste Values" would in the UI.
# This collapses the formula to its literal text result without touching
# any cell formatting/style.
$rng = $ws.Range("B2:E51")
$rng.Copy()
$rng.PasteSpecial(-4163)
$excel.CutCopyMode = 0
